$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "282.12"
Set-TextValue "D3" "20.59"
Set-TextValue "D4" "6.253"
Set-TextValue "D5" "0.06142"
Set-TextValue "D6" "3.577"
Set-TextValue "D7" "6.568"
Set-TextValue "D8" "1.501"
Set-TextValue "D9" "0.8189"
Set-TextValue "D10" "0.01382"
Set-TextValue "D11" "0.1634"
Set-TextValue "D12" "0.08400"
Set-TextValue "D14" "0.03181"
Set-TextValue "D15" "0.09131"
Set-TextValue "D16" "3.709"
Set-TextValue "D17" "0.001641"
Set-TextValue "D18" "0.04708"
Set-TextValue "D19" "0.006414"
Set-TextValue "D20" "0.006157"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.769"
Set-TextValue "D25" "0.3355"
Set-TextValue "D40" "0.04684"
Set-TextValue "D41" "0.007205"

# Row 42/43: swap CEJI <-> BKEXToken entries, with updated data
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1100"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003512"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D45" "0.00006614"
Set-TextValue "D48" "0.002947"
